$d = $word.ActiveDocument

# Title paragraph: "Fort Washington Ave Rehab Overview"
#               -> "Fort Washington Ave Rehab Consolidation Overview"
$d.Content.Find.Execute("Fort Washington Ave Rehab Overview", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Fort Washington Ave Rehab Consolidation Overview", 2)

# Body paragraph: "The Fort Washington Avenue Rehab consist of five developments in the "
#              -> "The Fort Washington Avenue Rehab Consolidation consist of five developments in the "
$d.Content.Find.Execute("The Fort Washington Avenue Rehab consist of five developments in the ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "The Fort Washington Avenue Rehab Consolidation consist of five developments in the ", 2)
